$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update username and password for rows 3 and 4
$ws.Range("C3").Value = "NTVNDM31"
$ws.Range("D3").Value = "1234@Welcome"
$ws.Range("C4").Value = "NTVNDM31"
$ws.Range("D4").Value = "1234@Welcome"

# Update the selected/active view of the sheet
$ws.Activate()
$ws.Range("C4:D4").Select()
